$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price value looks numeric (e.g. "1.00", "0.624") need the
# column pre-formatted as Text, otherwise Excel auto-converts the assigned
# string into a number/loses the trailing zeros - matching real Excel COM
# behaviour for Range.Value assignment.
$textFormatCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D14", "D18", "D19", "D20", "D21", "D22", "D24", "D25", "D26", "D27", "D28", "D30", "D31", "D33", "D34", "D35", "D36", "D37", "D38", "D40", "D41", "D42", "D43", "D44", "D45", "D47", "D48", "D50", "D51")
foreach ($cellRef in $textFormatCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '66.651.26'
$ws.Range("E2").Value = '  -4.21%  '
$ws.Range("D3").Value = '3.338.06'
$ws.Range("E3").Value = '  -1.71%  '
$ws.Range("D5").Value = '574.04'
$ws.Range("E5").Value = '  -3.36%  '
$ws.Range("D6").Value = '181.02'
$ws.Range("E6").Value = '  -5.61%  '
$ws.Range("D7").Value = '0.624'
$ws.Range("E7").Value = '  +2.79%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").Value = '0.130'
$ws.Range("E9").Value = '  -3.51%  '
$ws.Range("D10").Value = '6.66'
$ws.Range("E10").Value = '  -1.84%  '
$ws.Range("D11").Value = '0.402'
$ws.Range("E11").Value = '  -4.13%  '
$ws.Range("D12").Value = '3.917.43'
$ws.Range("E12").Value = '  -1.69%  '
$ws.Range("E13").Value = '  -0.58%  '
$ws.Range("D14").Value = '27.07'
$ws.Range("E14").Value = '  -6.04%  '
$ws.Range("D15").Value = '66.747.65'
$ws.Range("E15").Value = '  -4.03%  '
$ws.Range("E16").Value = '  -2.42%  '
$ws.Range("D17").Value = '3.354.60'
$ws.Range("E17").Value = '  -2.28%  '
$ws.Range("D18").Value = '438.34'
$ws.Range("E18").Value = '  -2.79%  '
$ws.Range("D19").Value = '5.69'
$ws.Range("E19").Value = '  -2.81%  '
$ws.Range("D20").Value = '13.56'
$ws.Range("E20").Value = '  -2.06%  '
$ws.Range("D21").Value = '7.61'
$ws.Range("E21").Value = '  -2.94%  '
$ws.Range("D22").Value = '73.58'
$ws.Range("E22").Value = '  -3.12%  '
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("D24").Value = '0.518'
$ws.Range("E24").Value = '  -1.00%  '
$ws.Range("D25").Value = '0.0000118'
$ws.Range("E25").Value = '  -4.28%  '
$ws.Range("D26").Value = '0.191'
$ws.Range("E26").Value = '  +0.28%  '
$ws.Range("D27").Value = '9.05'
$ws.Range("E27").Value = '  -4.65%  '
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  -0.33%  '
$ws.Range("E29").Value = '  -2.24%  '
$ws.Range("D30").Value = '22.86'
$ws.Range("E30").Value = '  -2.79%  '
$ws.Range("D31").Value = '5.31'
$ws.Range("E31").Value = '  -6.33%  '
$ws.Range("E32").Value = '  +0.03%  '
$ws.Range("D33").Value = '6.81'
$ws.Range("E33").Value = '  -2.84%  '
$ws.Range("D34").Value = '1.23'
$ws.Range("E34").Value = '  -4.70%  '
$ws.Range("D35").Value = '162.55'
$ws.Range("E35").Value = '  -1.69%  '
$ws.Range("D36").Value = '1.49'
$ws.Range("E36").Value = '  -5.88%  '
$ws.Range("D37").Value = '27.49'
$ws.Range("E37").Value = '  -2.45%  '
$ws.Range("D38").Value = '1.84'
$ws.Range("E38").Value = '  -5.83%  '
$ws.Range("D39").Value = '2.818.57'
$ws.Range("E39").Value = '  +2.22%  '
$ws.Range("D40").Value = '0.794'
$ws.Range("E40").Value = '  -2.96%  '
$ws.Range("D41").Value = '4.43'
$ws.Range("E41").Value = '  -4.14%  '
$ws.Range("D42").Value = '6.22'
$ws.Range("E42").Value = '  -5.80%  '
$ws.Range("D43").Value = '40.19'
$ws.Range("E43").Value = '  -2.47%  '
$ws.Range("D44").Value = '0.0667'
$ws.Range("E44").Value = '  -3.49%  '
$ws.Range("D45").Value = '24.46'
$ws.Range("E45").Value = '  -4.67%  '
$ws.Range("E46").Value = '  -6.96%  '
$ws.Range("D47").Value = '320.81'
$ws.Range("E47").Value = '  -6.22%  '
$ws.Range("D48").Value = '0.0274'
$ws.Range("E48").Value = '  -3.87%  '
$ws.Range("E49").Value = '  +0.48%  '
$ws.Range("D50").Value = '0.980'
$ws.Range("E50").Value = '  -4.33%  '
$ws.Range("D51").Value = '6.17'
$ws.Range("E51").Value = '  -3.09%  '
